$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4996
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("H58").Value = 7619.5
$ws.Range("I58").Value = 300
$ws.Range("K58").Value = 900
$ws.Range("M58").Value = -750
$ws.Range("H76").Value = 3749.75
$ws.Range("I76").Value = 3666.3333
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 3666.3333
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -3351.3333
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 3749.75
$ws.Range("I79").Value = 3666.3333
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 3666.3333
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -2574.3333
$ws.Range("N79").Value = -6184
$ws.Range("H86").Value = 3986
$ws.Range("I86").Value = 3996.5
$ws.Range("J86").Value = 3979
$ws.Range("K86").Value = 3996.5
$ws.Range("L86").Value = 3979
$ws.Range("M86").Value = -2873.5
$ws.Range("N86").Value = -6225
$ws.Range("H89").Value = 3986
$ws.Range("I89").Value = 3996.5
$ws.Range("J89").Value = 3979
$ws.Range("K89").Value = 19982.5
$ws.Range("L89").Value = 19895
$ws.Range("M89").Value = -14366.5
$ws.Range("N89").Value = -31127
$ws.Range("H113").Value = 1430288.2
$ws.Range("H137").Value = 1867.4193
$ws.Range("J137").Value = 2486
$ws.Range("L137").Value = 7458
$ws.Range("N137").Value = -12558
$ws.Range("H138").Value = 3673.422
$ws.Range("J138").Value = 3751.513
$ws.Range("L138").Value = 11254.539
$ws.Range("N138").Value = -21534.539
$ws.Range("H141").Value = 7467
$ws.Range("I141").Value = 7256
$ws.Range("J141").Value = 9999
$ws.Range("K141").Value = 21768
$ws.Range("L141").Value = 29997
$ws.Range("M141").Value = -16588
$ws.Range("N141").Value = -40357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4515.4307
$ws.Range("I32").Value = 4571.922
$ws.Range("K32").Value = 4571.922
$ws.Range("M32").Value = -4284.922
$ws.Range("H61").Value = 2326.182
$ws.Range("I61").Value = 1366.6666
$ws.Range("J61").Value = 2686
$ws.Range("K61").Value = 1366.6666
$ws.Range("L61").Value = 2686
$ws.Range("M61").Value = -1154.6666
$ws.Range("N61").Value = -3110
$ws.Range("H102").Value = 4833
$ws.Range("I102").Value = 4833
$ws.Range("K102").Value = 4833
$ws.Range("M102").Value = -3211
$ws.Range("H132").Value = 2762.8333
$ws.Range("I132").Value = 2220.8
$ws.Range("J132").Value = 3666.2222
$ws.Range("K132").Value = 6662.400000000001
$ws.Range("L132").Value = 10998.6666
$ws.Range("M132").Value = -4132.400000000001
$ws.Range("N132").Value = -16058.6666
$ws.Range("H136").Value = 2326.182
$ws.Range("I136").Value = 1366.6666
$ws.Range("J136").Value = 2686
$ws.Range("K136").Value = 4099.9998
$ws.Range("L136").Value = 8058
$ws.Range("M136").Value = -1549.9998
$ws.Range("N136").Value = -13158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 757.6667
$ws.Range("I22").Value = 769.4
$ws.Range("K22").Value = 769.4
$ws.Range("M22").Value = -596.4
$ws.Range("H64").Value = 1239.75
$ws.Range("J64").Value = 1319.6666
$ws.Range("L64").Value = 1319.6666
$ws.Range("N64").Value = -1769.6666
$ws.Range("H67").Value = 1239.75
$ws.Range("J67").Value = 1319.6666
$ws.Range("L67").Value = 1319.6666
$ws.Range("N67").Value = -2879.6666
$ws.Range("H86").Value = 2193.6
$ws.Range("I86").Value = 1993.1428
$ws.Range("K86").Value = 1993.1428
$ws.Range("M86").Value = -870.1428000000001
$ws.Range("H89").Value = 2193.6
$ws.Range("I89").Value = 1993.1428
$ws.Range("K89").Value = 9965.714
$ws.Range("M89").Value = -4349.714
$ws.Range("H94").Value = 1448.5555
$ws.Range("I94").Value = 910.1429000000001
$ws.Range("K94").Value = 910.1429000000001
$ws.Range("M94").Value = -459.1429000000001
$ws.Range("H134").Value = 3004.3333
$ws.Range("I134").Value = 3010.2856
$ws.Range("J134").Value = 2996
$ws.Range("K134").Value = 9030.856800000001
$ws.Range("L134").Value = 8988
$ws.Range("M134").Value = -6495.856800000001
$ws.Range("N134").Value = -14058

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1486
$ws.Range("I62").Value = 1486
$ws.Range("K62").Value = 1486
$ws.Range("M62").Value = -862
$ws.Range("H65").Value = 1486
$ws.Range("I65").Value = 1486
$ws.Range("K65").Value = 7430
$ws.Range("M65").Value = -4310
$ws.Range("H86").Value = 8226.904
$ws.Range("J86").Value = 6844.5454
$ws.Range("L86").Value = 6844.5454
$ws.Range("N86").Value = -9090.545399999999
$ws.Range("H89").Value = 8226.904
$ws.Range("J89").Value = 6844.5454
$ws.Range("L89").Value = 34222.727
$ws.Range("N89").Value = -45454.727
$ws.Range("H107").Value = 1982
$ws.Range("J107").Value = 2553.4167
$ws.Range("L107").Value = 2553.4167
$ws.Range("N107").Value = -6393.4167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1993.6666
$ws.Range("J75").Value = 2340.5
$ws.Range("L75").Value = 7021.5
$ws.Range("N75").Value = -9017.5
$ws.Range("H78").Value = 1993.6666
$ws.Range("J78").Value = 2340.5
$ws.Range("L78").Value = 21064.5
$ws.Range("N78").Value = -31048.5
$ws.Range("H80").Value = 171499
$ws.Range("I80").Value = 4996
$ws.Range("K80").Value = 14988
$ws.Range("M80").Value = -14052
$ws.Range("H83").Value = 171499
$ws.Range("I83").Value = 4996
$ws.Range("K83").Value = 44964
$ws.Range("M83").Value = -40284
$ws.Range("H103").Value = 1219.2222
$ws.Range("I103").Value = 724.5
$ws.Range("J103").Value = 1360.5714
$ws.Range("K103").Value = 2173.5
$ws.Range("L103").Value = 4081.7142
$ws.Range("M103").Value = -1294.5
$ws.Range("N103").Value = -5839.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5399.8
$ws.Range("I70").Value = 4999.75
$ws.Range("K70").Value = 4999.75
$ws.Range("M70").Value = -4729.75
$ws.Range("H73").Value = 5399.8
$ws.Range("I73").Value = 4999.75
$ws.Range("K73").Value = 4999.75
$ws.Range("M73").Value = -4063.75
$ws.Range("H132").Value = 4091.889
$ws.Range("J132").Value = 4191.154
$ws.Range("L132").Value = 12573.462
$ws.Range("N132").Value = -17633.462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6360
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 6360
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 6360
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6950
$ws.Range("H27").Value = 6360
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 6360
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 6360
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -6574
$ws.Range("H82").Value = 1385.4
$ws.Range("I82").Value = 1475.3334
$ws.Range("K82").Value = 1475.3334
$ws.Range("M82").Value = -1114.3334
$ws.Range("H85").Value = 1385.4
$ws.Range("I85").Value = 1475.3334
$ws.Range("K85").Value = 1475.3334
$ws.Range("M85").Value = -227.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 23334.334
$ws.Range("I4").Value = 21001.2
$ws.Range("J4").Value = 35000
$ws.Range("K4").Value = 21001.2
$ws.Range("L4").Value = 35000
$ws.Range("M4").Value = -20888.2
$ws.Range("N4").Value = -35226
